# "resultados e legendas novas"
# Adds a new "legend"/results block (columns Q:R) mirroring the existing
# L:M block, fills in several missing "Folha" (H column) percentages, and
# adjusts the view/column widths to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Folha" (leaf) severity values in column H ---------------------
$ws.Range("H10").Value = 6.19
$ws.Range("H14").Value = 9.06
$ws.Range("H18").Value = 6.33
$ws.Range("H22").Value = 11.89
$ws.Range("H26").Value = 3
$ws.Range("H30").Value = 12.37
$ws.Range("H34").Value = 3.14

# --- New legend block in columns Q:R (mirrors L:M around rows 11-12) ----
$ws.Range("Q11").Formula = "=SUM(H32:H33)"
$ws.Range("R11").Value = 100
$ws.Range("Q12").Value = 1.05
$ws.Range("R12").Formula = "=(R11*Q12)/Q11"
$ws.Range("R12").NumberFormat = "#,##0.00"

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 6.17
$ws.Columns.Item(18).ColumnWidth = 8.67

# --- View: scroll back to top-left, move selection to H35 ----------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H35").Select() | Out-Null
